# Reto 3 - Entrega final
# Fill in the remaining measurement data on the "Datos" sheet: the X-axis
# (number of requests: 1, 5, 10) for the M1-Memoria mini table, and the
# full M2-Tiempo / M2-Memoria mini tables, which had been left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# --- Table16 (Maquina 1 - Memoria), A6:G9 -> fill the "Maquina 1" (x) column
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 10
$ws.Range("A7:A9").Font.Bold = $false

# --- Table17 (Maquina 2 - Tiempo), A14:G17 -> fill all data rows
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = 61.734
$ws.Range("C15").Value = 10.861
$ws.Range("D15").Value = 13.452
$ws.Range("E15").Value = 611.75
$ws.Range("F15").Value = 7750.09
$ws.Range("G15").Value = 38441.611

$ws.Range("A16").Value = 5
$ws.Range("B16").Value = 242.182
$ws.Range("C16").Value = 133.849
$ws.Range("D16").Value = 94.811
$ws.Range("E16").Value = 4082.07
$ws.Range("F16").Value = 30949.489
$ws.Range("G16").Value = 511123.856

$ws.Range("A17").Value = 10
$ws.Range("B17").Value = 414.271
$ws.Range("C17").Value = 277.961
$ws.Range("D17").Value = 130.059
$ws.Range("E17").Value = 4805.525
$ws.Range("F17").Value = 132028.608
$ws.Range("G17").Value = 1637760.284

$ws.Range("B15:G17").NumberFormat = "#,##0"

# Grow Table17 by one row (matches the new Table173 extent used downstream)
$lo17 = $ws.ListObjects.Item("Table17")
$lo17.Resize($ws.Range("A14:G18"))
$lo17.Name = "Table173"

# --- Table18 (Maquina 2 - Memoria), A19:G22 -> fill all data rows
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = 21.875
$ws.Range("C20").Value = 30.672
$ws.Range("D20").Value = 14.852
$ws.Range("E20").Value = 2570.507
$ws.Range("F20").Value = 2651.764
$ws.Range("G20").Value = 299144.385

$ws.Range("A21").Value = 5
$ws.Range("B21").Value = 21.398
$ws.Range("C21").Value = 29.555
$ws.Range("D21").Value = 16.156
$ws.Range("E21").Value = 5001.804
$ws.Range("F21").Value = 3972.615
$ws.Range("G21").Value = 333723.485

$ws.Range("A22").Value = 10
$ws.Range("B22").Value = 21.305
$ws.Range("C22").Value = 31.07
$ws.Range("D22").Value = 16.68
$ws.Range("E22").Value = 5619.327
$ws.Range("F22").Value = 4155.744
$ws.Range("G22").Value = 6517063.856

$ws.Range("B20:G22").NumberFormat = "#,##0"

$lo18 = $ws.ListObjects.Item("Table18")
$lo18.Name = "Table184"

$wb.Save()
